# Update the weekly Time Record sheet with new findings:
#  - mark an item as "NA YET" and give it an actual-time value
#  - fill in the remaining daily "Time Worked" entries
#  - add weekly totals (SUM formulas) in column G
#  - move the active selection to reflect the new end of data entry

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Actual Time" values for the Modeling-related rows
$ws.Range("C5").Value = "NA YET"
$ws.Range("C6").Value = 2

# Fill in the remaining "Time Worked" values for the second week
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = 3
$ws.Range("F13").Value = 2
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = 4

# Weekly total formulas
$ws.Range("G8").Formula = "=SUM(F2:F8)"
$ws.Range("G15").Formula = "=SUM(F9:F15)"

# Reflect the new active cell/selection
[void]$ws.Range("G16").Select()
